# Add an "enabled" column to the Properties sheet (new optional-argument
# flag for each property row), and make the Properties sheet the active
# (selected) tab with C2 as the active cell - mirrors the upstream commit
# "added optional arguments to script".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Properties")

$ws.Range("C1").Value = "enabled"
$ws.Range("C2").Value = $true
$ws.Range("C3").Value = $true

$ws.Activate()
$ws.Range("C2").Select()
